{"js": "// Load the document body's paragraph collection so we can locate the two\n// trailing empty paragraphs that sit right before the final section break.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The document ends with two completely empty paragraphs, right after the\n// paragraph that contains the \"tal como lo menciona GitHub\" screenshot.\n// Re-use those two existing (empty) paragraphs for the two new sentences\n// instead of inserting brand-new ones, so no extra paragraphs are created.\nconst count = paragraphs.items.length;\nconst firstPara = paragraphs.items[count - 2];\nconst secondPara = paragraphs.items[count - 1];\n\nfirstPara.insertText(\n  \"Dicho esto. Se realiza la muestra del c\u00f3digo fuente a continuaci\u00f3n, para que se pueda apreciar cada l\u00ednea.\",\n  Word.InsertLocation.replace\n);\nsecondPara.insertText(\n  \"De m\u00e1s est\u00e1 decir que el c\u00f3digo est\u00e1 en su totalidad comentado, excepto los m\u00e9todos con eventos o funciones f\u00e1cil de entender a primera vista para quien conoce sobre l\u00f3gica de programaci\u00f3n o el lenguaje de C#.\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\n// After those two paragraphs, append 13 plain empty paragraphs.\nlet lastPara = secondPara;\nfor (let i = 0; i < 13; i++) {\n  lastPara = lastPara.insertParagraph(\"\", Word.InsertLocation.after);\n}\nawait context.sync();\n\n// Finally, append 5 more empty paragraphs that only carry a negative right\n// indent (w:ind w:right=\"-427\", i.e. -427 twips = -21.35 points).\nfor (let i = 0; i < 5; i++) {\n  lastPara = lastPara.insertParagraph(\"\", Word.InsertLocation.after);\n  lastPara.rightIndent = -21.35;\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The document ends with two completely empty paragraphs, right after the\n# paragraph that contains the \"tal como lo menciona GitHub\" screenshot.\n# Re-use those two existing (empty) paragraphs for the two new sentences\n# instead of inserting brand-new ones, so no extra paragraphs are created.\n$count = $d.Paragraphs.Count\n$firstPara = $d.Paragraphs($count - 1)\n$secondPara = $d.Paragraphs($count)\n\n$firstPara.Range.Text = \"Dicho esto. Se realiza la muestra del c\u00f3digo fuente a continuaci\u00f3n, para que se pueda apreciar cada l\u00ednea.\"\n$secondPara.Range.Text = \"De m\u00e1s est\u00e1 decir que el c\u00f3digo est\u00e1 en su totalidad comentado, excepto los m\u00e9todos con eventos o funciones f\u00e1cil de entender a primera vista para quien conoce sobre l\u00f3gica de programaci\u00f3n o el lenguaje de C#.\"\n\n# After those two paragraphs, append 13 plain empty paragraphs.\nfor ($i = 0; $i -lt 13; $i++) {\n    $n = $d.Paragraphs.Count\n    $lastPara = $d.Paragraphs($n)\n    $lastPara.Range.InsertParagraphAfter()\n}\n\n# Finally, append 5 more empty paragraphs that only carry a negative right\n# indent (w:ind w:right=\"-427\", i.e. -427 twips = -21.35 points).\nfor ($i = 0; $i -lt 5; $i++) {\n    $n = $d.Paragraphs.Count\n    $lastPara = $d.Paragraphs($n)\n    $lastPara.Range.InsertParagraphAfter()\n    $n2 = $d.Paragraphs.Count\n    $newPara = $d.Paragraphs($n2)\n    $newPara.Range.ParagraphFormat.RightIndent = -21.35\n}\n"}
